$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H11").Value = 32.285713
$ws.Range("I11").Value = 32.285713
$ws.Range("K11").Value = 32.285713
$ws.Range("M11").Value = 107.714287
$ws.Range("H33").Value = 2667.6667
$ws.Range("I33").Value = 1004
$ws.Range("K33").Value = 1004
$ws.Range("M33").Value = -775
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = $null
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 2423.2856
$ws.Range("I32").Value = 1252.7
$ws.Range("K32").Value = 1252.7
$ws.Range("M32").Value = -965.7
$ws.Range("H132").Value = 1865
$ws.Range("I132").Value = 1865
$ws.Range("K132").Value = 5595
$ws.Range("M132").Value = -3065
$ws = $wb.Worksheets.Item(3)
$ws.Range("H105").Value = 15471.667
$ws.Range("I105").Value = 15471.667
$ws.Range("K105").Value = 15471.667
$ws.Range("M105").Value = -13724.667
$ws.Range("H107").Value = 2263.625
$ws.Range("J107").Value = 1487.25
$ws.Range("L107").Value = 1487.25
$ws.Range("N107").Value = -5327.25
$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 72.14286
$ws.Range("I7").Value = 73.333336
$ws.Range("J7").Value = 71.25
$ws.Range("K7").Value = 73.333336
$ws.Range("L7").Value = 71.25
$ws.Range("M7").Value = 39.666664
$ws.Range("N7").Value = -297.25
$ws.Range("H58").Value = 3999
$ws.Range("I58").Value = 1331.6666
$ws.Range("K58").Value = 1331.6666
$ws.Range("M58").Value = -1128.6666
$ws.Range("H136").Value = 3999
$ws.Range("I136").Value = 1331.6666
$ws.Range("K136").Value = 3994.9998
$ws.Range("M136").Value = -1444.9998
$ws = $wb.Worksheets.Item(5)
$ws.Range("H5").Value = 535
$ws.Range("J5").Value = 802
$ws.Range("L5").Value = 2406
$ws.Range("N5").Value = -2630
$ws.Range("H26").Value = 149.5
$ws.Range("J26").Value = 149.5
$ws.Range("L26").Value = 448.5
$ws.Range("N26").Value = -1024.5
$ws.Range("H68").Value = 308.75
$ws.Range("I68").Value = 300
$ws.Range("J68").Value = 335
$ws.Range("K68").Value = 900
$ws.Range("L68").Value = 1005
$ws.Range("M68").Value = -89
$ws.Range("N68").Value = -2627
$ws.Range("H71").Value = 308.75
$ws.Range("I71").Value = 300
$ws.Range("J71").Value = 335
$ws.Range("K71").Value = 2700
$ws.Range("L71").Value = 3015
$ws.Range("M71").Value = 1356
$ws.Range("N71").Value = -11127
$ws.Range("H109").Value = 30
$ws.Range("J109").Value = 30
$ws.Range("L109").Value = 90
$ws.Range("N109").Value = -2170
$ws.Range("H112").Value = 2000
$ws.Range("J112").Value = 2000
$ws.Range("L112").Value = 6000
$ws.Range("N112").Value = -8216
$ws.Range("H123").Value = 1332.6666
$ws.Range("I123").Value = 1332.6666
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 3997.9998
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = $null
$ws.Range("N123").Value = -1547.9998
$ws.Range("H135").Value = 535
$ws.Range("J135").Value = 802
$ws.Range("L135").Value = 7218
$ws.Range("N135").Value = -12288
$ws.Range("H140").Value = 1347
$ws.Range("I140").Value = 1347
$ws.Range("K140").Value = 4041
$ws.Range("M140").Value = 1139
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 5900
$ws.Range("I70").Value = 5900
$ws.Range("K70").Value = 5900
$ws.Range("M70").Value = -5630
$ws.Range("H73").Value = 5900
$ws.Range("I73").Value = 5900
$ws.Range("K73").Value = 5900
$ws.Range("M73").Value = -4964
$ws.Range("H122").Value = 1999.5
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = $null
$ws = $wb.Worksheets.Item(7)
$ws.Range("H4").Value = 12600
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = 16400
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 16400
$ws.Range("M4").Value = -4887
$ws.Range("N4").Value = -16626
$ws.Range("H7").Value = 11000.8
$ws.Range("I7").Value = 3004
$ws.Range("J7").Value = 13000
$ws.Range("K7").Value = 3004
$ws.Range("L7").Value = 13000
$ws.Range("M7").Value = -2892
$ws.Range("N7").Value = -13224
$ws.Range("H28").Value = 12600
$ws.Range("I28").Value = 5000
$ws.Range("J28").Value = 16400
$ws.Range("K28").Value = 5000
$ws.Range("L28").Value = 16400
$ws.Range("M28").Value = -4768
$ws.Range("N28").Value = -16864
$ws.Range("H31").Value = 26000
$ws.Range("I31").Value = 4000
$ws.Range("J31").Value = 48000
$ws.Range("K31").Value = 4000
$ws.Range("L31").Value = 48000
$ws.Range("M31").Value = -3752
$ws.Range("N31").Value = -48496
$ws.Range("H37").Value = 12600
$ws.Range("I37").Value = 5000
$ws.Range("J37").Value = 16400
$ws.Range("K37").Value = 5000
$ws.Range("L37").Value = 16400
$ws.Range("M37").Value = -4893
$ws.Range("N37").Value = -16614
$ws.Range("H126").Value = 11000.8
$ws.Range("I126").Value = 3004
$ws.Range("J126").Value = 13000
$ws.Range("K126").Value = 9012
$ws.Range("L126").Value = 39000
$ws.Range("M126").Value = -6542
$ws.Range("N126").Value = -43940
$ws.Range("H136").Value = 9756.416999999999
$ws.Range("I136").Value = 5452.3335
$ws.Range("K136").Value = 16357.0005
$ws.Range("M136").Value = -13807.0005
$ws = $wb.Worksheets.Item(8)
$ws.Range("H28").Value = 150000
$ws.Range("J28").Value = 150000
$ws.Range("L28").Value = 150000
$ws.Range("N28").Value = -150696
$ws.Range("H122").Value = 790.875
$ws.Range("I122").Value = 754.5
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 2263.5
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = 186.5
$ws.Range("N122").Value = -7600
$ws.Range("H132").Value = 3235.875
$ws.Range("I132").Value = 2981.1667
$ws.Range("K132").Value = 8943.500100000001
$ws.Range("M132").Value = -6413.500100000001
$ws.Range("H136").Value = 3270.8333
$ws.Range("I136").Value = 2280.75
$ws.Range("J136").Value = 5251
$ws.Range("K136").Value = 6842.25
$ws.Range("L136").Value = 15753
$ws.Range("M136").Value = -4292.25
$ws.Range("N136").Value = -20853
